# Refresh the crypto price/volume table (rows 2-51) with the latest
# coinranking.com snapshot values, matching the GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.011.16"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -6.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.544.31"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.62"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.43"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -6.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -3.39%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -5.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.68"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -7.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -4.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.65"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -4.79%  "
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.937.10"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.583.59"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.868"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -4.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.07"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -4.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.060.50"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -6.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.06"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0975"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.59"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.64"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.83"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -10.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.92"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -5.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.98"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -5.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.48"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("E30").Value = "  -5.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.95"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.30"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.15"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.36"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -6.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0793"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -5.27%  "
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.83"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +7.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.22"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +8.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0310"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -4.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.88"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.076.72"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "84.57"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -10.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.89"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.59"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.795.37"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.23"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.68"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -3.96%  "
